# Update the "想去人数" (interested-people count) figures that were
# refreshed when the gh-pages data was regenerated.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 205
$ws1.Range("F7").Value  = 1052
$ws1.Range("F14").Value = 585
$ws1.Range("F18").Value = 1259
$ws1.Range("F20").Value = 2888
$ws1.Range("F22").Value = 707
$ws1.Range("F28").Value = 3137

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 205
$ws4.Range("F13").Value = 1052
$ws4.Range("F25").Value = 585
$ws4.Range("F29").Value = 1259
$ws4.Range("F31").Value = 2888
$ws4.Range("F33").Value = 707
$ws4.Range("F41").Value = 3137
